$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Near the end of the document: drop the paragraph that duplicated the
#    bold title text ("Play Champion of the Underworld for Free -
#    Review and Analysis") and update the italic paragraph that follows
#    it with the new image-prompt text.
#
#    Done first (while paragraph indices are still the original ones) to
#    keep the lookups simple; iterate paragraphs directly (rather than
#    Content.Find) and skip paragraph 1 so the real Heading1 title at
#    the top of the document is never touched.
# ---------------------------------------------------------------------

$oldTitleText = "Play Champion of the Underworld for Free - Review and Analysis"
$dupIndex = -1
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd("`r")
    if ($t -eq $oldTitleText) {
        $dupIndex = $i
    }
}
if ($dupIndex -gt 0) {
    $d.Paragraphs.Item($dupIndex).Range.Delete()
}

$oldDescription = "Looking for a unique and enjoyable gaming experience? Read our review on Champion of the Underworld and play for free today."
$newDescription = "Please create a bold and eye-catching feature image for Champion of the Underworld that portrays the unique theme and graphics of the game. The image should be in a cartoon style and feature a happy Maya warrior sporting glasses, as this is one of the key symbols in the game. The warrior could be shown battling Hades or one of the other monstrous creatures that appear in the game's underworld setting. The background of the image could be a dark and mysterious underworld filled with glowing crystals and other fantastical elements that represent the game's theme. Make sure the colors are vibrant and appealing to help attract players and convey the exciting gameplay of this Yggdrasil slot."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd("`r")
    if ($t -eq $oldDescription) {
        $descRange = $d.Range($p.Range.Start, $p.Range.End - 1)
        $descRange.Text = $newDescription
    }
}

# ---------------------------------------------------------------------
# 2) Insert a new "Meta description" paragraph right after the Heading1
#    title paragraph.
#
#    We want the new paragraph to end up with NO explicit paragraph
#    properties (matching a plain body / Normal paragraph, just like the
#    diff shows). Splitting an existing paragraph (InsertParagraphAfter/
#    Before) always copies the pPr of the paragraph being split, so
#    instead we clone an already-"plain" paragraph (the first body
#    paragraph under "Gameplay Mechanics" - no w:pPr at all) via
#    Copy/Paste, which keeps its (lack of) paragraph formatting, and
#    then overwrite its runs with the exact run structure we need.
# ---------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)

$plainParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -ne "" -and $p.Style.NameLocal -eq "Normal" -and $plainParaIndex -eq -1) {
        $plainParaIndex = $i
    }
}
$plainPara = $d.Paragraphs.Item($plainParaIndex)
$plainPara.Range.Copy()

$insertPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$insertPoint.Paste()

$newPara = $d.Paragraphs.Item(2)
$newParaText = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Looking for a unique and enjoyable gaming experience? Read our review on Champion of the Underworld and play for free today.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newParaText.InsertXML($metaXml)
